# finish dev of upgradeMilitaryTech
# Adds two new columns (C: INT_techPointPerLevel, D: STR_building) to the
# militaryTechs sheet, fills them in for every tech row, and nudges the
# window/selection view state to match the author's last-saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header/data-cell formatting (fill/border/alignment) from
# column B onto the new columns C and D so the new cells pick up the same
# cellXfs styles Excel would have carried over.
$ws.Range("B1:B17").Copy()
$ws.Range("C1:D17").PasteSpecial(-4122)  # xlPasteFormats

# --- column D (STR_building) -------------------------------------------
# Each military tech name is "<attacker>_<defender>"; the building that
# trains/researches it is derived from the attacker side.
$ws.Range("D1").Value = "STR_building"

$buildingByPrefix = @{
    "infantry" = "trainingGround"
    "archer"   = "hunterHall"
    "cavalry"  = "stable"
    "siege"    = "workshop"
}

for ($row = 2; $row -le 17; $row++) {
    $techName = $ws.Range("A$row").Value2
    $prefix = $techName.Split("_")[0]
    $building = $buildingByPrefix[$prefix]
    $ws.Cells.Item($row, 4).Value = $building
}

# --- column C (INT_techPointPerLevel) -----------------------------------
$ws.Range("C1").Value = "INT_techPointPerLevel"

for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 3).Value = 50
}

# --- view state -----------------------------------------------------------
[void]$ws.Range("C2").Select()

$wb.Windows.Item(1).Left = 7960
$wb.Windows.Item(1).Top = 3920
